$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 2.997429241610044)
    3  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 9.295990156953671)
    4  = @(0.127881588408715, 0.3127903958511391, 26.21740644021617, 8.660232485948974, 0, 35.318310910425)
    5  = @(0.3048080303191223, 0.04240448674262143, 3.900430680208489, 8.660232485948974, 0, 12.90787568321921)
    6  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 4.429675500412797)
    7  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    8  = @(1.459612070389937, 10.29869402782916, 0.1575252929769615, 645.3272768299601, 0, 657.2431082211561)
    9  = @(0.127881588408715, 1.667794583268128, 3.900430680208489, 645.3272768299601, 0, 651.0233836818454)
    10 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 9.295990156953671)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
